# Generate Report for Handoff
#
# Re-orders the three tracked-file rows (by file name) on every sheet and
# refreshes the "d265b39b..." row's handoff/handback status, matching a
# new localization-status report run. Row order becomes:
#   row2 -> ffffe21ea4ca-75ec-4446-a428-3736196fd5e6.md
#   row3 -> ffffff1b38a086-09e9-4be9-b34b-f75de2b996e7.md
#   row4 -> d265b39b-0772-403e-b473-d76686770375.md   (now "Ready for handoff")
#
# Hyperlinked cells keep their original r:id (so their link target stays
# put) while the cell text / hyperlink display text is updated in place -
# this mirrors how the source report-generator tool re-wrote the sheet.

$wb = $excel.ActiveWorkbook

function Set-CellText {
    param($ws, [string]$addr, [string]$value)

    $ws.Range($addr).Value = $value

    # If this cell carries a hyperlink, keep the hyperlink (and its r:id /
    # target) but refresh the displayed text to match the new cell value.
    $col = $addr -replace '[0-9]+$', ''
    $row = $addr -replace '^[A-Za-z]+', ''
    $target = '$' + $col + '$' + $row

    $hls = @($ws.Hyperlinks)
    foreach ($h in $hls) {
        if ($h.Range.Address() -eq $target) {
            $h.TextToDisplay = $value
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellText $wsOverview "A2" "ffffe21ea4ca-75ec-4446-a428-3736196fd5e6.md"
Set-CellText $wsOverview "B2" "e2e\ffffe21ea4ca-75ec-4446-a428-3736196fd5e6.md"
Set-CellText $wsOverview "G2" "2016-08-24 11:05:35"

Set-CellText $wsOverview "A3" "ffffff1b38a086-09e9-4be9-b34b-f75de2b996e7.md"
Set-CellText $wsOverview "B3" "e2e\ffffff1b38a086-09e9-4be9-b34b-f75de2b996e7.md"

Set-CellText $wsOverview "A4" "d265b39b-0772-403e-b473-d76686770375.md"
Set-CellText $wsOverview "B4" "e2e\d265b39b-0772-403e-b473-d76686770375.md"
Set-CellText $wsOverview "E4" "Ready for handoff"
Set-CellText $wsOverview "F4" "Ready for handoff"
Set-CellText $wsOverview "G4" "2016-08-24 11:08:21"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-CellText $wsZh "A2" "ffffe21ea4ca-75ec-4446-a428-3736196fd5e6.md"
Set-CellText $wsZh "G2" "30d6b8d0-65a7-4bca-8fda-8549ca162525.922bb0427acb37c44ea70b4f17270d17034ef070.zh-cn.xlf"
Set-CellText $wsZh "H2" "2016-08-24 11:05:30"
Set-CellText $wsZh "I2" "30d6b8d0-65a7-4bca-8fda-8549ca162525.md"
Set-CellText $wsZh "J2" "30d6b8d0-65a7-4bca-8fda-8549ca162525.922bb0427acb37c44ea70b4f17270d17034ef070.zh-cn.xlf"
Set-CellText $wsZh "K2" "2016-08-24 11:05:47"

Set-CellText $wsZh "A3" "ffffff1b38a086-09e9-4be9-b34b-f75de2b996e7.md"
Set-CellText $wsZh "F3" "True"

Set-CellText $wsZh "A4" "d265b39b-0772-403e-b473-d76686770375.md"
Set-CellText $wsZh "C4" "Ready for handoff"
Set-CellText $wsZh "F4" "False"
Set-CellText $wsZh "G4" "d265b39b-0772-403e-b473-d76686770375.58e823c7542409bfd107ae584079abe28d687504.zh-cn.xlf"
Set-CellText $wsZh "H4" "2016-08-24 11:08:16"
Set-CellText $wsZh "I4" "d265b39b-0772-403e-b473-d76686770375.md"
Set-CellText $wsZh "J4" "d265b39b-0772-403e-b473-d76686770375.58e823c7542409bfd107ae584079abe28d687504.zh-cn.xlf"
Set-CellText $wsZh "K4" "2016-08-24 11:07:38"
Set-CellText $wsZh "P4" "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0dd174a7b11733fc71afb87d3fe12db9d52b3fce/e2e/d265b39b-0772-403e-b473-d76686770375.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc8891b031880b72f8b3277957e26be10e64a247/e2e/d265b39b-0772-403e-b473-d76686770375.md."

# Error Detail column needs to be wider to show the new message.
$wsZh.Range("P1").ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-CellText $wsDe "A2" "ffffe21ea4ca-75ec-4446-a428-3736196fd5e6.md"
Set-CellText $wsDe "G2" "30d6b8d0-65a7-4bca-8fda-8549ca162525.922bb0427acb37c44ea70b4f17270d17034ef070.de-de.xlf"
Set-CellText $wsDe "H2" "2016-08-24 11:05:35"
Set-CellText $wsDe "I2" "30d6b8d0-65a7-4bca-8fda-8549ca162525.md"
Set-CellText $wsDe "J2" "30d6b8d0-65a7-4bca-8fda-8549ca162525.922bb0427acb37c44ea70b4f17270d17034ef070.de-de.xlf"
Set-CellText $wsDe "K2" "2016-08-24 11:05:54"

Set-CellText $wsDe "A3" "ffffff1b38a086-09e9-4be9-b34b-f75de2b996e7.md"
Set-CellText $wsDe "F3" "True"

Set-CellText $wsDe "A4" "d265b39b-0772-403e-b473-d76686770375.md"
Set-CellText $wsDe "C4" "Ready for handoff"
Set-CellText $wsDe "F4" "False"
Set-CellText $wsDe "G4" "d265b39b-0772-403e-b473-d76686770375.58e823c7542409bfd107ae584079abe28d687504.de-de.xlf"
Set-CellText $wsDe "H4" "2016-08-24 11:08:21"
Set-CellText $wsDe "I4" "d265b39b-0772-403e-b473-d76686770375.md"
Set-CellText $wsDe "J4" "d265b39b-0772-403e-b473-d76686770375.58e823c7542409bfd107ae584079abe28d687504.de-de.xlf"
Set-CellText $wsDe "K4" "2016-08-24 11:07:46"
Set-CellText $wsDe "P4" "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0dd174a7b11733fc71afb87d3fe12db9d52b3fce/e2e/d265b39b-0772-403e-b473-d76686770375.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc8891b031880b72f8b3277957e26be10e64a247/e2e/d265b39b-0772-403e-b473-d76686770375.md."

# Error Detail column needs to be wider to show the new message.
$wsDe.Range("P1").ColumnWidth = 39.166666666666664
